$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Status text: "In Translation" -> "Handed back: in sync with en-US"
#    This shared string is referenced from the Overview sheet (E2/F2/E3/F3)
#    as well as the "Status" column (C2/C3) on both the zh-cn and de-de
#    sheets, so update every one of those cells to the new text.
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Record the handback for each language/file pair.
#    For every row: J = link to the handed-back markdown file (same link as
#    column A), K = name of the generated handback xliff, L = handback
#    timestamp.
# ---------------------------------------------------------------------------
$mdUrlB6 = "https://github.com/OpenLocalizationTestOrg/ol-test1/blob/173e366b951b1b06c6c32c5ccc7973559891d33b/e2e/b6be4a1d-a99d-49b4-bd1f-c6684c1736b6.md"
$mdUrlC2 = "https://github.com/OpenLocalizationTestOrg/ol-test1/blob/173e366b951b1b06c6c32c5ccc7973559891d33b/e2e/c288ce95-522f-4983-9a83-23e0e7a12296.md"
$mdNameB6 = "b6be4a1d-a99d-49b4-bd1f-c6684c1736b6.md"
$mdNameC2 = "c288ce95-522f-4983-9a83-23e0e7a12296.md"

# zh-cn
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("J2"), $mdUrlB6, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdNameB6) | Out-Null
$wsZhCn.Range("K2").Value = "b6be4a1d-a99d-49b4-bd1f-c6684c1736b6.cda0b3f32037fd9f1ac93e9aeda6cec7f33f3c80.zh-cn.xlf"
$wsZhCn.Range("L2").Value = "2017-01-03 04:31:52"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("J3"), $mdUrlC2, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdNameC2) | Out-Null
$wsZhCn.Range("K3").Value = "c288ce95-522f-4983-9a83-23e0e7a12296.c3d59441392122d3234ca06da4966323e96f1b8b.zh-cn.xlf"
$wsZhCn.Range("L3").Value = "2017-01-03 04:31:52"

# de-de
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("J2"), $mdUrlB6, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdNameB6) | Out-Null
$wsDeDe.Range("K2").Value = "b6be4a1d-a99d-49b4-bd1f-c6684c1736b6.cda0b3f32037fd9f1ac93e9aeda6cec7f33f3c80.de-de.xlf"
$wsDeDe.Range("L2").Value = "2017-01-03 04:32:04"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("J3"), $mdUrlC2, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdNameC2) | Out-Null
$wsDeDe.Range("K3").Value = "c288ce95-522f-4983-9a83-23e0e7a12296.c3d59441392122d3234ca06da4966323e96f1b8b.de-de.xlf"
$wsDeDe.Range("L3").Value = "2017-01-03 04:32:04"

# Give the newly populated "Latest Target File" cells the same Hyperlink
# look-and-feel as the existing source-file links in column A.
$wsZhCn.Range("J2").Style = "HyperLink"
$wsZhCn.Range("J3").Style = "HyperLink"
$wsDeDe.Range("J2").Style = "HyperLink"
$wsDeDe.Range("J3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# 3. Widen the columns that now hold the handback info / longer status text
#    so the new content is readable.
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.1   # E
$wsOverview.Columns.Item(6).ColumnWidth = 29.1   # F

$wsZhCn.Columns.Item(3).ColumnWidth = 29.1    # C - Status
$wsZhCn.Columns.Item(10).ColumnWidth = 39.1   # J - Latest Target File
$wsZhCn.Columns.Item(11).ColumnWidth = 39.1   # K - Latest Handback File

$wsDeDe.Columns.Item(3).ColumnWidth = 29.1    # C - Status
$wsDeDe.Columns.Item(10).ColumnWidth = 39.1   # J - Latest Target File
$wsDeDe.Columns.Item(11).ColumnWidth = 39.1   # K - Latest Handback File
